# Apply latest Katalon test-run timestamps (Date column) and refreshed
# Pass/Fail Result for the VerifyReceiptPageLabelsCredit sheet, row 2,
# reflecting the most recent execution after the Windows 11 update.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("PayNowCC")
$ws.Range("B2").Value = "Mon Aug 04 23:52:19 IST 2025"
$ws.Range("B3").Value = "Mon Aug 04 23:53:04 IST 2025"
$ws.Range("B4").Value = "Mon Aug 04 23:53:51 IST 2025"
$ws.Range("B5").Value = "Mon Aug 04 23:54:39 IST 2025"
$ws.Range("B6").Value = "Mon Aug 04 23:55:25 IST 2025"
$ws.Range("B7").Value = "Mon Aug 04 23:56:11 IST 2025"
$ws.Range("B8").Value = "Mon Aug 04 23:56:58 IST 2025"
$ws.Range("B9").Value = "Mon Aug 04 23:57:43 IST 2025"

$ws = $wb.Worksheets.Item("VerifyConfirmPageLabelsCredit")
$ws.Range("B2").Value = "Tue Aug 05 00:20:17 IST 2025"
$ws.Range("B3").Value = "Tue Aug 05 00:20:58 IST 2025"

$ws = $wb.Worksheets.Item("VerifyReceiptPageLabelsCredit")
$ws.Range("A2").Value = "Pass"
$ws.Range("B2").Value = "Tue Aug 05 22:07:21 IST 2025"
$ws.Range("B3").Value = "Tue Aug 05 22:03:33 IST 2025"

$ws = $wb.Worksheets.Item("VerifyConfirmPageLabelsPC")
$ws.Range("B2").Value = "Tue Aug 05 00:23:09 IST 2025"
$ws.Range("B3").Value = "Tue Aug 05 00:23:55 IST 2025"

$ws = $wb.Worksheets.Item("VerifyReceiptPageLabelsPC")
$ws.Range("B2").Value = "Tue Aug 05 00:44:04 IST 2025"
$ws.Range("B3").Value = "Tue Aug 05 00:45:02 IST 2025"

$ws = $wb.Worksheets.Item("PayNowCorpSCF")
$ws.Range("B2").Value = "Mon Aug 04 23:21:04 IST 2025"
$ws.Range("B3").Value = "Mon Aug 04 23:22:07 IST 2025"
$ws.Range("B4").Value = "Mon Aug 04 23:23:13 IST 2025"
$ws.Range("B5").Value = "Mon Aug 04 23:24:21 IST 2025"

$ws = $wb.Worksheets.Item("PayNowCreditSCF")
$ws.Range("B2").Value = "Mon Aug 04 23:30:32 IST 2025"
$ws.Range("B3").Value = "Mon Aug 04 23:31:35 IST 2025"
$ws.Range("B4").Value = "Mon Aug 04 23:32:39 IST 2025"
$ws.Range("B5").Value = "Mon Aug 04 23:33:45 IST 2025"

$ws = $wb.Worksheets.Item("PayNowPersonalCheckSCF")
$ws.Range("B2").Value = "Mon Aug 04 23:39:17 IST 2025"
$ws.Range("B3").Value = "Mon Aug 04 23:40:20 IST 2025"
$ws.Range("B4").Value = "Mon Aug 04 23:41:28 IST 2025"
$ws.Range("B5").Value = "Mon Aug 04 23:42:32 IST 2025"

$ws = $wb.Worksheets.Item("PayNowPersonalSavingsSCF")
$ws.Range("B2").Value = "Mon Aug 04 23:43:45 IST 2025"
$ws.Range("B3").Value = "Mon Aug 04 23:44:49 IST 2025"
$ws.Range("B4").Value = "Mon Aug 04 23:45:53 IST 2025"
$ws.Range("B5").Value = "Mon Aug 04 23:46:56 IST 2025"

$ws = $wb.Worksheets.Item("PayNowCreditDCF")
$ws.Range("B2").Value = "Mon Aug 04 23:26:09 IST 2025"
$ws.Range("B3").Value = "Mon Aug 04 23:27:14 IST 2025"
$ws.Range("B4").Value = "Mon Aug 04 23:28:21 IST 2025"
$ws.Range("B5").Value = "Mon Aug 04 23:29:23 IST 2025"

$ws = $wb.Worksheets.Item("PayNowCorpDCF")
$ws.Range("B2").Value = "Mon Aug 04 23:16:48 IST 2025"
$ws.Range("B3").Value = "Mon Aug 04 23:17:55 IST 2025"
$ws.Range("B4").Value = "Mon Aug 04 23:18:57 IST 2025"
$ws.Range("B5").Value = "Mon Aug 04 23:20:01 IST 2025"

$ws = $wb.Worksheets.Item("PayNowPC")
$ws.Range("B2").Value = "Mon Aug 04 23:10:39 IST 2025"

$ws = $wb.Worksheets.Item("PayNowPersonalCheckDCF")
$ws.Range("B2").Value = "Mon Aug 04 23:48:05 IST 2025"
$ws.Range("B3").Value = "Mon Aug 04 23:49:10 IST 2025"
$ws.Range("B4").Value = "Mon Aug 04 23:50:12 IST 2025"
$ws.Range("B5").Value = "Mon Aug 04 23:51:14 IST 2025"

$ws = $wb.Worksheets.Item("MaxAmountErrorCC")
$ws.Range("B3").Value = "Tue Aug 05 01:26:07 IST 2025"

$ws = $wb.Worksheets.Item("MaxAmountErrorCorp")
$ws.Range("B2").Value = "Tue Aug 05 01:26:57 IST 2025"
$ws.Range("B3").Value = "Tue Aug 05 01:27:45 IST 2025"

$ws = $wb.Worksheets.Item("MaxAmountErrorPC")
$ws.Range("B2").Value = "Tue Aug 05 01:28:40 IST 2025"
$ws.Range("B3").Value = "Tue Aug 05 01:29:26 IST 2025"

$ws = $wb.Worksheets.Item("MaxAmountErrorPS")
$ws.Range("B2").Value = "Tue Aug 05 01:30:11 IST 2025"
$ws.Range("B3").Value = "Tue Aug 05 01:31:01 IST 2025"

$ws = $wb.Worksheets.Item("MinAmountErrorPC")
$ws.Range("B2").Value = "Tue Aug 05 01:34:58 IST 2025"
$ws.Range("B3").Value = "Tue Aug 05 01:35:47 IST 2025"

$ws = $wb.Worksheets.Item("MinAmountErrorCC")
$ws.Range("B2").Value = "Tue Aug 05 01:31:47 IST 2025"
$ws.Range("B3").Value = "Tue Aug 05 01:32:36 IST 2025"

$ws = $wb.Worksheets.Item("MinAmountErrorCorp")
$ws.Range("B2").Value = "Tue Aug 05 01:33:21 IST 2025"
$ws.Range("B3").Value = "Tue Aug 05 01:34:11 IST 2025"

$ws = $wb.Worksheets.Item("MinAmountErrorPS")
$ws.Range("B2").Value = "Tue Aug 05 01:36:36 IST 2025"
$ws.Range("B3").Value = "Tue Aug 05 01:37:22 IST 2025"

$ws = $wb.Worksheets.Item("PayNowPS")
$ws.Range("B2").Value = "Mon Aug 04 23:11:45 IST 2025"
$ws.Range("B3").Value = "Mon Aug 04 23:12:42 IST 2025"
$ws.Range("B4").Value = "Mon Aug 04 23:13:36 IST 2025"
$ws.Range("B5").Value = "Mon Aug 04 23:14:25 IST 2025"
$ws.Range("B6").Value = "Mon Aug 04 23:15:13 IST 2025"
$ws.Range("B7").Value = "Mon Aug 04 23:15:58 IST 2025"

$ws = $wb.Worksheets.Item("OverAndUnderPayCredit")
$ws.Range("B2").Value = "Tue Aug 05 00:02:02 IST 2025"
$ws.Range("B3").Value = "Tue Aug 05 00:02:54 IST 2025"
$ws.Range("B4").Value = "Tue Aug 05 00:03:38 IST 2025"
$ws.Range("B5").Value = "Tue Aug 05 00:04:22 IST 2025"

$ws = $wb.Worksheets.Item("OverAndUnderPayPC")
$ws.Range("B2").Value = "Tue Aug 05 00:05:09 IST 2025"
$ws.Range("B3").Value = "Tue Aug 05 00:05:53 IST 2025"
$ws.Range("B4").Value = "Tue Aug 05 00:06:43 IST 2025"
$ws.Range("B5").Value = "Tue Aug 05 00:07:27 IST 2025"

$ws = $wb.Worksheets.Item("OverAndUnderPayPS")
$ws.Range("B2").Value = "Tue Aug 05 00:08:21 IST 2025"
$ws.Range("B3").Value = "Tue Aug 05 00:09:19 IST 2025"
$ws.Range("B4").Value = "Tue Aug 05 00:10:08 IST 2025"
$ws.Range("B5").Value = "Tue Aug 05 00:10:55 IST 2025"

$ws = $wb.Worksheets.Item("OverAndUnderPayCorp")
$ws.Range("B2").Value = "Mon Aug 04 23:58:32 IST 2025"
$ws.Range("B3").Value = "Mon Aug 04 23:59:25 IST 2025"
$ws.Range("B4").Value = "Tue Aug 05 00:00:17 IST 2025"
$ws.Range("B5").Value = "Tue Aug 05 00:01:11 IST 2025"

$ws = $wb.Worksheets.Item("NoModifyAmountCorp")
$ws.Range("B2").Value = "Wed Aug 06 00:31:54 IST 2025"
$ws.Range("B3").Value = "Wed Aug 06 00:32:58 IST 2025"

$ws = $wb.Worksheets.Item("NoModifyAmountPC")
$ws.Range("B2").Value = "Wed Aug 06 00:34:35 IST 2025"
$ws.Range("B3").Value = "Wed Aug 06 00:35:35 IST 2025"

$ws = $wb.Worksheets.Item("NoModifyAmountPS")
$ws.Range("B2").Value = "Wed Aug 06 00:37:00 IST 2025"
$ws.Range("B3").Value = "Wed Aug 06 00:38:00 IST 2025"

$ws = $wb.Worksheets.Item("NoModifyAmountCC")
$ws.Range("B2").Value = "Tue Aug 05 01:49:23 IST 2025"
$ws.Range("B3").Value = "Tue Aug 05 01:50:11 IST 2025"

$ws = $wb.Worksheets.Item("NoOverPayErrorCC")
$ws.Range("B2").Value = "Tue Aug 05 01:39:44 IST 2025"
$ws.Range("B3").Value = "Tue Aug 05 01:40:34 IST 2025"

$ws = $wb.Worksheets.Item("NoOverPayErrorPC")
$ws.Range("B2").Value = "Tue Aug 05 01:43:02 IST 2025"
$ws.Range("B3").Value = "Tue Aug 05 01:43:47 IST 2025"

$ws = $wb.Worksheets.Item("PayNowCorp")
$ws.Range("B2").Value = "Mon Aug 04 23:25:23 IST 2025"

$ws = $wb.Worksheets.Item("NoOverPayErrorCorp")
$ws.Range("B2").Value = "Tue Aug 05 01:41:28 IST 2025"
$ws.Range("B3").Value = "Tue Aug 05 01:42:15 IST 2025"

$ws = $wb.Worksheets.Item("NoOverPayErrorPS")
$ws.Range("B2").Value = "Tue Aug 05 01:44:36 IST 2025"
$ws.Range("B3").Value = "Tue Aug 05 01:45:21 IST 2025"

$ws = $wb.Worksheets.Item("NoUnderPayErrorPS")
$ws.Range("B2").Value = "Tue Aug 05 01:48:36 IST 2025"

$ws = $wb.Worksheets.Item("NoUnderPayErrorPC")
$ws.Range("B2").Value = "Tue Aug 05 01:47:42 IST 2025"

$ws = $wb.Worksheets.Item("NoUnderPayErrorCC")
$ws.Range("B2").Value = "Tue Aug 05 01:46:07 IST 2025"

$ws = $wb.Worksheets.Item("NoUnderPayErrorCorp")
$ws.Range("B2").Value = "Tue Aug 05 01:46:57 IST 2025"

$ws = $wb.Worksheets.Item("CardExpiredErrorCC")
$ws.Range("B2").Value = "Tue Aug 05 01:22:45 IST 2025"
$ws.Range("B3").Value = "Tue Aug 05 01:23:34 IST 2025"

$ws = $wb.Worksheets.Item("CardNotAcceptedErrorCC")
$ws.Range("B2").Value = "Tue Aug 05 01:24:28 IST 2025"
$ws.Range("B3").Value = "Tue Aug 05 01:25:18 IST 2025"

$ws = $wb.Worksheets.Item("MRFCorp")
$ws.Range("B2").Value = "Tue Aug 05 01:38:09 IST 2025"
$ws.Range("B3").Value = "Tue Aug 05 01:39:03 IST 2025"

$ws = $wb.Worksheets.Item("VerifyConfirmPageLabelsCorp")
$ws.Range("B2").Value = "Tue Aug 05 00:21:44 IST 2025"
$ws.Range("B3").Value = "Tue Aug 05 00:22:28 IST 2025"

$ws = $wb.Worksheets.Item("VerifyReceiptPageLabelsCorp")
$ws.Range("B2").Value = "Tue Aug 05 22:10:24 IST 2025"
$ws.Range("B3").Value = "Tue Aug 05 22:11:29 IST 2025"

$ws = $wb.Worksheets.Item("VerifyConfirmPageLabelsPS")
$ws.Range("B2").Value = "Tue Aug 05 00:24:49 IST 2025"
$ws.Range("B3").Value = "Tue Aug 05 00:25:35 IST 2025"

$ws = $wb.Worksheets.Item("VerifyReceiptPageLabelsPS")
$ws.Range("B2").Value = "Tue Aug 05 00:46:01 IST 2025"
$ws.Range("B3").Value = "Tue Aug 05 00:47:01 IST 2025"
